$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Simple text swaps (paragraphs without adjacent empty runs) ---

$d.Content.Find.Execute("Play Eternal Shogi Free: Unique Gameplay & Bonus Features", $true, $false, $false, $false, $false, $true, 1, $false, "Play Eternal Shogi Free - Unique Gameplay and Exciting Bonus Features", 2)

$d.Content.Find.Execute("Discover Eternal Shogi, a unique slot game with bonus features based on traditional board games. Play for free and enjoy suitable betting ranges for all levels.", $true, $false, $false, $false, $false, $true, 1, $false, "Discover the unique gameplay and exciting bonus features in Eternal Shogi. Play for free now!", 2)

# --- Bullet paragraphs: each has a leading empty <w:r/> that a plain text
# replace on the run would otherwise merge away. Rebuild the whole
# paragraph (excluding its paragraph mark) via InsertXML so the empty run
# survives untouched, matching the source structure exactly. ---

function Set-BulletParagraph($findText, $newText) {
    foreach ($p in $d.Paragraphs) {
        $ptext = $p.Range.Text.TrimEnd([char]13, [char]10, [char]7)
        if ($ptext -eq $findText) {
            $xml = "<w:p $wns><w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>$newText</w:t></w:r></w:p>"
            $p.Range.InsertXML($xml)
        }
    }
}

Set-BulletParagraph "Suitable betting ranges for both low and high rollers" "Wide range of betting options for all players"
Set-BulletParagraph "Symbols and bonus features inspired by traditional board games" "Engaging symbols and bonus features"
Set-BulletParagraph "Limited to 10 paylines" "Limited number of paylines"
